# Update the carjacking-by-neighborhood-by-month workbook with the
# 2022-03-28 data refresh (data "through March 20" instead of "March 19").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-03-20"

# Update the column header text (shared string) for the "March 2022" column.
$ws.Range("B1").Value = "March 2022 (through March 20)"

# --- Updated counts for existing cells ---
$ws.Range("E3").Value = 5     # Austin, March 2021
$ws.Range("E4").Value = 9     # North Lawndale, March 2021
$ws.Range("B5").Value = 6     # Garfield Park, March 2022
$ws.Range("N9").Value = 3     # Chicago Lawn, March 2018
$ws.Range("E11").Value = 2    # Englewood, March 2021
$ws.Range("B15").Value = 4    # Humboldt Park, March 2022
$ws.Range("E26").Value = 3    # Grand Crossing, March 2021
$ws.Range("K33").Value = 2    # Belmont Cragin, March 2019
$ws.Range("H49").Value = 3    # Little Village, March 2020

# --- Newly populated cells (previously blank) ---
$ws.Range("N20").Value = 1    # West Lawn, March 2018
$ws.Range("B24").Value = 1    # Wicker Park, March 2022
$ws.Range("B27").Value = 1    # Calumet Heights, March 2022
$ws.Range("Q34").Value = 1    # River North, March 2017
$ws.Range("T49").Value = 1    # Little Village, March 2016
